$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (testCase04) - fill in the Actual and Status columns
# Actual mirrors the Expected value, Status is marked as Pass
$ws.Range("E5").Value = $ws.Range("D5").Text
$ws.Range("F5").Value = "Pass"
